$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.840.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.303.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.23"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.63%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.05%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.39"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0907"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.49"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.79%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.17%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.32"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.652.39"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.301.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.796.91"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.38"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +29.41%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.49"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.06"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.55"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.68"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.64%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.36"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.04"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +19.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.42"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.66"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.50"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0874"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.46%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.115"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.32%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0349"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.75"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.89%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.66%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.36"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.64%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.39"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.32"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.39%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.52"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.26%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.739.14"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.33"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.80"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.12"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.39%  "
